# Weekly data refresh: a new week's price observation is inserted at the
# top of the "Feria Lagunitas de Puerto Montt - Poroto verde" series
# (row 80), pushing the existing rows 80-92 down to 81-93. The sheet's
# used-range grows from A1:R92 to A1:R93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 80..92 down one row, opening up a fresh row 80.
$ws.Rows.Item(80).Insert()

# Populate the newly opened row 80 with this week's record.
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 44782
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = 100112031
$ws.Range("G80").Value = "Poroto verde"
$ws.Range("H80").Value = "Magnum"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 35
$ws.Range("K80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("M80").Value = 40000
$ws.Range("N80").Value = "`$/malla 25 kilos"
$ws.Range("O80").Value = "Perú"
$ws.Range("P80").Value = 1600
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"
